$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 30 de Junio de 2020 a las 10:58"
$ws.Range("B6").Value = 647849
$ws.Range("C6").Value = 6693
$ws.Range("D6").Value = 412650
$ws.Range("E6").Value = 225879
$ws.Range("G6").Value = 154
$ws.Range("H6").Value = 9320
$ws.Range("A20").Value = "Banglades"
$ws.Range("B20").Value = 145483
$ws.Range("C20").Value = 3682
$ws.Range("D20").Value = 59624
$ws.Range("E20").Value = 84012
$ws.Range("G20").Value = 64
$ws.Range("H20").Value = 1847
$ws.Range("A21").Value = "Sudafrica"
$ws.Range("B21").Value = 144264
$ws.Range("D21").Value = 70614
$ws.Range("E21").Value = 71121
$ws.Range("H21").Value = 2529
$ws.Range("B42").Value = 34393
$ws.Range("C42").Value = 239
$ws.Range("D42").Value = 21281
$ws.Range("E42").Value = 11649
$ws.Range("G42").Value = 19
$ws.Range("H42").Value = 1463
$ws.Range("B47").Value = 31359
$ws.Range("C47").Value = 121
$ws.Range("E47").Value = 16690
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 735
$ws.Range("B53").Value = 24688
$ws.Range("C53").Value = 247
$ws.Range("D53").Value = 17272
$ws.Range("E53").Value = 7096
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 320
$ws.Range("B72").Value = 8866
$ws.Range("C72").Value = 4
$ws.Range("E72").Value = 479
$ws.Range("B76").Value = 7214
$ws.Range("C76").Value = 5
$ws.Range("E76").Value = 286
$ws.Range("A79").Value = "El Salvador"
$ws.Range("B79").Value = 6438
$ws.Range("C79").Value = 265
$ws.Range("D79").Value = 3770
$ws.Range("E79").Value = 2494
$ws.Range("G79").Value = 10
$ws.Range("H79").Value = 174
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 6209
$ws.Range("D80").Value = 2427
$ws.Range("E80").Value = 3484
$ws.Range("H80").Value = 298
$ws.Range("A81").Value = "Kenia"
$ws.Range("B81").Value = 6190
$ws.Range("D81").Value = 2013
$ws.Range("E81").Value = 4033
$ws.Range("H81").Value = 144
$ws.Range("B118").Value = 1667
$ws.Range("C118").Value = 2
$ws.Range("E118").Value = 175
$ws.Range("B120").Value = 1600
$ws.Range("C120").Value = 15
$ws.Range("E120").Value = 105
$ws.Range("A168").Value = "Namibia"
$ws.Range("B168").Value = 203
$ws.Range("C168").Value = 7
$ws.Range("D168").Value = 24
$ws.Range("E168").Value = 179
$ws.Range("H168").Value = 0
$ws.Range("A169").Value = "Islas Caimanes"
$ws.Range("B169").Value = 199
$ws.Range("D169").Value = 189
$ws.Range("E169").Value = 9
$ws.Range("H169").Value = 1
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
